$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "59.462.15"
$ws.Range("E2").Value = "  +2.51%  "

# Row 3
$ws.Range("D3").Value = "3.014.19"
$ws.Range("E3").Value = "  +1.43%  "

# Row 4
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "565.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.34%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.85%  "

# Row 8
$ws.Range("E8").Value = "  +1.36%  "

# Row 9
$ws.Range("D9").Value = "3.000.98"
$ws.Range("E9").Value = "  +1.18%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.134"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.63%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.23"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +7.21%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.459"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.23%  "

# Row 13
$ws.Range("E13").Value = "  +3.87%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.63%  "

# Row 15
$ws.Range("E15").Value = "  +2.30%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.28%  "

# Row 17
$ws.Range("D17").Value = "3.504.06"
$ws.Range("E17").Value = "  +1.20%  "

# Row 18
$ws.Range("D18").Value = "3.009.45"
$ws.Range("E18").Value = "  +1.26%  "

# Row 19
$ws.Range("D19").Value = "59.427.42"
$ws.Range("E19").Value = "  +2.60%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "433.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.35%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.33%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.725"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.45%  "

# Row 23
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.64"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.95%  "

# Row 24
$ws.Range("B24").Value = "Uniswap"
$ws.Range("C24").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.39%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.83"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.46%  "

# Row 26
$ws.Range("E26").Value = "  +0.04%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +11.92%  "

# Row 28
$ws.Range("E28").Value = "  +0.05%  "

# Row 29
$ws.Range("E29").Value = "  +2.38%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.93"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.51%  "

# Row 31
$ws.Range("E31").Value = "  +2.23%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.15"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.06%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.100"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.24%  "

# Row 34
$ws.Range("B34").Value = "Mantle"
$ws.Range("C34").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.01"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.92%  "

# Row 35
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.01"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.15%  "

# Row 36
$ws.Range("D36").Value = "0.0₃0765"
$ws.Range("E36").Value = "  +10.03%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.13"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.51%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "49.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.58%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.70"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.27%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.77"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.23%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "410.30"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.62%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0356"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.26%  "

# Row 43
$ws.Range("D43").Value = "2.774.57"
$ws.Range("E43").Value = "  +3.25%  "

# Row 44
$ws.Range("E44").Value = "  -0.39%  "

# Row 45
$ws.Range("E45").Value = "  +4.68%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "35.43"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +23.58%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "123.90"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.46%  "

# Row 49
$ws.Range("B49").Value = "Fetch.AI"
$ws.Range("C49").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.03"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.95%  "

# Row 50
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.111"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.58%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.64"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.24%  "
